# Time delay and lag 2 step optimization
# Update the simulated/optimized numeric results in Sheet1 with refreshed
# values produced by a new optimization run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.3750999659704599
$ws.Range("C2").Value = 3.4500595774128766
$ws.Range("D2").Value = 1.6391627272293354
$ws.Range("E2").Value = 1.4912153409914097

$ws.Range("B3").Value = 1.839299632548925
$ws.Range("C3").Value = 0.57421817686101884
$ws.Range("D3").Value = 1.9248988078259972
$ws.Range("E3").Value = 1.4168087968476042

$ws.Range("B4").Value = 2.4111537307867583
$ws.Range("C4").Value = 0.38172367380654204
$ws.Range("D4").Value = 3.0028397097794528
$ws.Range("E4").Value = 1.662715776824724

$ws.Range("B5").Value = 0.38704723084047982
$ws.Range("C5").Value = 14.85420109889974
$ws.Range("D5").Value = 0.37570128233747402
$ws.Range("E5").Value = 0.49193177896086732

$ws.Range("E6").Value = 29.996407810101477

$ws.Range("B9").Value = 32.978970526181556
$ws.Range("C9").Value = 3.0000255431164389
$ws.Range("D9").Value = 32.974399796881222
$ws.Range("E9").Value = 32.981668565336498

$ws.Range("E10").Value = 3.0072943573859146

$ws.Range("B11").Value = 0.72125116029461611
$ws.Range("C11").Value = 1.333666917992703
$ws.Range("D11").Value = 0.90524422118296255
$ws.Range("E11").Value = 0.72475219688574244

$ws.Range("B12").Value = 3.9540006386741626
$ws.Range("C12").Value = 17.423704598910113
$ws.Range("D12").Value = 4.4077336711713668
$ws.Range("E12").Value = 5.2485168304383238

$ws.Range("B13").Value = 9.4438898024387541
$ws.Range("C13").Value = 29.524046030795393
$ws.Range("D13").Value = 5.5448721289041547
$ws.Range("E13").Value = 8.4511642667044917

$ws.Range("C14").Value = [double]"8.1127857221791728e-15"

$ws.Range("B15").Value = 1.0944388980243875
$ws.Range("C15").Value = 1.2952404603079539
$ws.Range("D15").Value = 1.0554488235124633
$ws.Range("E15").Value = 1.0845116426670449

$ws.Range("B16").Value = 1.9880341774986883
$ws.Range("C16").Value = 4.0849543957641687
$ws.Range("D16").Value = 2.6802684809369128
$ws.Range("E16").Value = 1.9759349973543752
